# cards_constructs.xlsx update: add "operation_X" construct row to the
# "construct definitions" sheet, per commit "updates cards_constructs with
# operationId".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("construct definitions")

# --- New row 16 plain-text columns -----------------------------------------
$ws.Range("A16").Value = "operation_X"
$ws.Range("B16").Value = "Simple"
$ws.Range("C16").Value = "X'th operation ID for the AnalysisMethod"
$ws.Range("D16").Value = "Analyses`$method_id ->  AnalysisMethods`$operation_id (X'th operation ID)"
$ws.Range("E16").Value = "opidXhere"

# Match the wrap-text style already used by columns C, D, F, G in the rest
# of the table (same as e.g. row 15).
$ws.Range("C16").WrapText = $true
$ws.Range("D16").WrapText = $true

# --- F16: example usage snippet (template), with the two fill-in-here ------
# placeholders highlighted in red, matching the convention used elsewhere in
# this workbook.
$f16 = @'
df3_analysisidhere <- df3_analysisidhere|>
           dplyr::filter(stat_name %in% c('n', 'p')) |>
           dplyr::mutate(operationid = dplyr::case_when(stat_name == 'n' ~ 'opid1here',
                                                                                         stat_name == 'p' ~ 'opid2here'))
'@
$ws.Range("F16").Value = $f16
$ws.Range("F16").WrapText = $true
$ws.Range("F16").Characters(176, 9).Font.Color = 255
$ws.Range("F16").Characters(297, 9).Font.Color = 255

# --- G16: example usage, resulting code populated with metadata ------------
$g16 = @'
df3_An_02 <- df3_An_02|>
                    dplyr::filter(stat_name %in% c('n', 'p')) |>
                    dplyr::mutate(operationid = dplyr::case_when(stat_name == 'n' ~ 'Mth_03_01_n',
                                                                                                   stat_name == 'p' ~ 'Mth_03_02_%'))
'@
$ws.Range("G16").Value = $g16
$ws.Range("G16").WrapText = $true
$ws.Range("G16").Characters(176, 11).Font.Color = 255
$ws.Range("G16").Characters(309, 11).Font.Color = 255

# Row 16 wraps to four lines, same height as other 4-line rows (e.g. row 3).
$ws.Rows.Item(16).RowHeight = 57.6

# --- Sheet-level cosmetics to match the saved view --------------------------
# Column F widened slightly (author grew it to fit the new content) and no
# longer auto-"best fit".
$ws.Columns.Item(6).ColumnWidth = 80.5

# Selection/scroll position saved along with the edit.
$ws.Activate()
$ws.Range("C6").Select()
